$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2: change Target cluster from "MuSCs" to "ECs" and refresh the numeric values
$ws.Range("D2").Value = "ECs"

$ws.Range("G2").Value = 1.524170333333333
$ws.Range("H2").Value = 4.572511
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03919966666666667
$ws.Range("N2").Value = 0.117599
$ws.Range("O2").Value = 0.2591884862911046
$ws.Range("P2").Value = 0.2591884862911046
$ws.Range("Q2").Value = 0.05974696900988889
$ws.Range("R2").Value = 0.537722721089
$ws.Range("S2").Value = 0.2591884862911046
$ws.Range("T2").Value = 0.2591884862911046

# Add new row 3 with the original "MuSCs" target cluster data, updated values
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pomc"
$ws.Range("C3").Value = "Mc4r"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.524170333333333
$ws.Range("H3").Value = 4.572511
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1120403333333333
$ws.Range("N3").Value = 0.336121
$ws.Range("O3").Value = 0.7408115137088953
$ws.Range("P3").Value = 0.7408115137088953
$ws.Range("Q3").Value = 0.1707685522034445
$ws.Range("R3").Value = 1.536916969831
$ws.Range("S3").Value = 0.7408115137088953
$ws.Range("T3").Value = 0.7408115137088953
